# Weekly update: insert the newest "Acelga" price record for
# "Vega Modelo de Temuco" right after the header block of existing
# records (it becomes the new row 188), pushing all the subsequent
# historical rows down by one. Excel's native row-insert semantics take
# care of shifting every existing row (188->189, ..., 215->216) and of
# growing the sheet's used range, so all that is left to do is stamp the
# brand-new row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 188:215 down to 189:216, leaving a blank row 188 behind.
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new weekly record.
$ws.Range("A188").Value = 10
$ws.Range("B188").Value = "Vega Modelo de Temuco"
$ws.Range("C188").Value = "La Araucanía"
$ws.Range("D188").Value = 44505
$ws.Range("E188").Value = 9
$ws.Range("F188").Value = 100112009
$ws.Range("G188").Value = "Acelga"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 65
$ws.Range("K188").Value = 6000
$ws.Range("L188").Value = 6000
$ws.Range("M188").Value = 6000
$ws.Range("N188").Value = "$/docena de atados (12 kilos)"
$ws.Range("O188").Value = "Provincia de Cautín"
$ws.Range("P188").Value = 500
$ws.Range("Q188").Value = 12
$ws.Range("R188").Value = "Hortaliza"
